# Rename user "lrezende" to "wlima" throughout the "grants por usuario" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")

# Column B (rows 59:112) holds the username that is referenced by the
# concatenation formulas in column D. Replacing the value(s) here updates
# every dependent formula once Excel recalculates.
$rng = $ws.Range("B59:B112")
$rng.Replace("lrezende", "wlima", 1) | Out-Null

$excel.CalculateFullRebuild()
